$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D, E and G hold numeric-looking values (prices,
# percentages, hour) that must stay plain text, matching the
# original inline-string cell contents -- force text format
# first so COM does not auto-coerce them into real numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '329.18'
$ws.Range("E2").Value = '0.67%'
$ws.Range("G2").Value = '22'

$ws.Range("D3").Value = '41.27'
$ws.Range("E3").Value = '2.30%'
$ws.Range("G3").Value = '22'

$ws.Range("D4").Value = '5.704'
$ws.Range("E4").Value = '-1.14%'
$ws.Range("G4").Value = '22'

$ws.Range("D5").Value = '0.08066'
$ws.Range("E5").Value = '0.68%'
$ws.Range("G5").Value = '22'

$ws.Range("D6").Value = '2.050'
$ws.Range("E6").Value = '6.58%'
$ws.Range("G6").Value = '22'

$ws.Range("D7").Value = '8.714'
$ws.Range("E7").Value = '0.00%'
$ws.Range("G7").Value = '22'

$ws.Range("D8").Value = '4.517'
$ws.Range("E8").Value = '-1.42%'
$ws.Range("G8").Value = '22'

$ws.Range("D9").Value = '2.952'
$ws.Range("E9").Value = '0.36%'
$ws.Range("G9").Value = '22'

$ws.Range("D10").Value = '0.9214'
$ws.Range("E10").Value = '-2.53%'
$ws.Range("G10").Value = '22'

$ws.Range("D11").Value = '0.1258'
$ws.Range("E11").Value = '0.59%'
$ws.Range("G11").Value = '22'

$ws.Range("D12").Value = '0.1942'
$ws.Range("E12").Value = '-0.74%'
$ws.Range("G12").Value = '22'

$ws.Range("D13").Value = '8.227'
$ws.Range("E13").Value = '-8.22%'
$ws.Range("G13").Value = '22'

$ws.Range("D14").Value = '0.09270'
$ws.Range("E14").Value = '1.14%'
$ws.Range("G14").Value = '22'

$ws.Range("D15").Value = '0.03669'
$ws.Range("E15").Value = '4.46%'
$ws.Range("G15").Value = '22'

$ws.Range("D16").Value = '0.1054'
$ws.Range("E16").Value = '10.02%'
$ws.Range("G16").Value = '22'

$ws.Range("D17").Value = '0.001300'
$ws.Range("E17").Value = '-0.12%'
$ws.Range("G17").Value = '22'

$ws.Range("D18").Value = '0.006320'
$ws.Range("E18").Value = '-1.69%'
$ws.Range("G18").Value = '22'

$ws.Range("D19").Value = '3.382'
$ws.Range("E19").Value = '0.50%'
$ws.Range("G19").Value = '22'

$ws.Range("E20").Value = '-1.49%'
$ws.Range("G20").Value = '22'

$ws.Range("D21").Value = '0.1417'
$ws.Range("E21").Value = '0.79%'
$ws.Range("G21").Value = '22'

$ws.Range("D22").Value = '0.2650'
$ws.Range("E22").Value = '9.44%'
$ws.Range("G22").Value = '22'

$ws.Range("D23").Value = '0.04432'
$ws.Range("E23").Value = '0.25%'
$ws.Range("G23").Value = '22'

$ws.Range("D24").Value = '0.001259'
$ws.Range("E24").Value = '-0.46%'
$ws.Range("G24").Value = '22'

$ws.Range("D25").Value = '0.004327'
$ws.Range("E25").Value = '0.34%'
$ws.Range("G25").Value = '22'

$ws.Range("E26").Value = '8.04%'
$ws.Range("G26").Value = '22'

$ws.Range("G27").Value = '22'

$ws.Range("G28").Value = '22'

$ws.Range("G29").Value = '22'

$ws.Range("G30").Value = '22'

$ws.Range("G31").Value = '22'

$ws.Range("G32").Value = '22'

$ws.Range("G33").Value = '22'

$ws.Range("G34").Value = '22'

$ws.Range("G35").Value = '22'

$ws.Range("G36").Value = '22'

$ws.Range("G37").Value = '22'

$ws.Range("G38").Value = '22'

$ws.Range("D39").Value = '0.02851'
$ws.Range("E39").Value = '17.95%'
$ws.Range("G39").Value = '22'

$ws.Range("E40").Value = '4.51%'
$ws.Range("G40").Value = '22'

$ws.Range("D41").Value = '0.007593'
$ws.Range("E41").Value = '1.24%'
$ws.Range("G41").Value = '22'

$ws.Range("D42").Value = '0.009933'
$ws.Range("E42").Value = '14.27%'
$ws.Range("G42").Value = '22'

$ws.Range("E43").Value = '0.02%'
$ws.Range("G43").Value = '22'

$ws.Range("D44").Value = '0.002111'
$ws.Range("E44").Value = '-0.23%'
$ws.Range("G44").Value = '22'

$ws.Range("D45").Value = '0.01180'
$ws.Range("E45").Value = '7.50%'
$ws.Range("G45").Value = '22'

$ws.Range("D46").Value = '0.00006728'
$ws.Range("E46").Value = '-2.87%'
$ws.Range("G46").Value = '22'

$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").Value = '-0.75%'
$ws.Range("G47").Value = '22'

$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '0.002279'
$ws.Range("E48").Value = '59.71%'
$ws.Range("G48").Value = '22'

$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '0.002992'
$ws.Range("E49").Value = '-5.73%'
$ws.Range("G49").Value = '22'

$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.75%'
$ws.Range("G50").Value = '22'

$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.75%'
$ws.Range("G51").Value = '22'
